# Generate Report for Handoff
#
# The localization CI tool re-ran: the translation finished, so the status
# text for the zh-cn / de-de locales flips from "In Translation" to
# "Ready for handoff", and the associated "generated" timestamps advance.
# The Status / per-locale-status columns are auto-fit to their contents by
# the report generator, so they widen to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Width (in "characters", the ColumnWidth unit) that best-fits the new,
# longer status string "Ready for handoff" for these report columns.
$statusColWidth = 16.3333333333333

# --- Overview sheet -------------------------------------------------
# E2 ("zh-cn" column) / F2 ("de-de" column): per-locale status
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-17 22:57:04"

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ------------------------------------------------------
# C2: Status
$wsZhCn.Range("C2").Value = "Ready for handoff"
# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-17 22:56:57"

$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet --------------------------------------------------------
# C2: Status
$wsDeDe.Range("C2").Value = "Ready for handoff"
# H2: Latest Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-08-17 22:57:04"

$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
